# Applies the "Added switcher part, updated nums" commit:
# adds a second JPEG-vs-AVIF comparison table (rows 29-39, cols B-H)
# below the existing data on Лист1, and nudges the window/view state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- New "JPEG vs AVIF" block -------------------------------------------

$ws.Range("B29").Value = "JPEG vs AVIF"

$ws.Range("B30").Value = "JPEG"
$ws.Range("F30").Value = "AVIF"

$ws.Range("B31").Value = "FAST 4G"
$ws.Range("C31").Value = "Double 4G"
$ws.Range("D31").Value = "No Trottling"
$ws.Range("F31").Value = "FAST 4G"
$ws.Range("G31").Value = "Double 4G"
$ws.Range("H31").Value = "No Trottling"

# JPEG measurements (B:D) and AVIF measurements (F:H), rows 32-37
$jpeg = @(
    @(811.2,   694.08, 146.58000000000001),
    @(792.15,  692.51, 144.51),
    @(820.08,  704.87, 157.13),
    @(782.8,   643.79999999999995, 159.35),
    @(771.88,  690.88, 149.24),
    @(820.84,  708.29, 154.12)
)
$avif = @(
    @(753.2,   694.38, 173.08),
    @(725.18,  694.53, 162.81),
    @(722.84,  675.82, 223.72),
    @(748.55,  715.88, 219),
    @(742.96,  673.42, 193.85),
    @(726.33,  683.84, 217.85)
)

for ($i = 0; $i -lt 6; $i++) {
    $row = 32 + $i
    $ws.Cells.Item($row, 2).Value = $jpeg[$i][0]
    $ws.Cells.Item($row, 3).Value = $jpeg[$i][1]
    $ws.Cells.Item($row, 4).Value = $jpeg[$i][2]

    $ws.Cells.Item($row, 6).Value = $avif[$i][0]
    $ws.Cells.Item($row, 7).Value = $avif[$i][1]
    $ws.Cells.Item($row, 8).Value = $avif[$i][2]
}

# ---- Averages (row 38) and standard deviations (row 39) -----------------
# B is a standalone formula; C:H (skipping the blank separator column E)
# is entered as one multi-cell range so Excel records it as a single shared
# formula, same pattern already used by the first table's row 25/26.

$ws.Range("A38").Value = "Среднее"
$ws.Range("B38").Formula = "=AVERAGE(B32:B37)"
$ws.Range("C38:H38").Formula = "=AVERAGE(C32:C37)"
$ws.Range("E38").ClearContents()

$ws.Range("A39").Value = "Стандартное отклонение"
$ws.Range("B39").Formula = "=STDEV(B32:B37)"
$ws.Range("C39:H39").Formula = "=STDEV(C32:C37)"
$ws.Range("E39").ClearContents()

# Match the numeric style (2 decimal places) used by the rest of the sheet.
$ws.Range("B38:D39").NumberFormat = "0.00"
$ws.Range("F38:H39").NumberFormat = "0.00"
$ws.Range("E38").NumberFormat = "0.00"
$ws.Range("E39").NumberFormat = "0.00"

# ---- Page setup (printable, A4-ish single sheet like the source) --------

$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# ---- View state: scroll to the new block, select the new total cell -----

$ws.Application.ActiveWindow.ScrollRow = 10
$ws.Range("F38").Select()
